# Applies the changes described by the commit "fix SyntheticDataPipeline and regenerate"
# to the synthetic_data workbook.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. License sheet: BSD title -> "BSD License"
# ---------------------------------------------------------------------------
$wsLicense = $wb.Worksheets.Item("License")
$wsLicense.Range("C2").Value = "BSD License"

# ---------------------------------------------------------------------------
# 2. RightsStatement sheet: "note" column (E) values got reshuffled between
#    rows (rows 2,4,5,6,7,8,10,11,12,13). Row 9 is untouched.
# ---------------------------------------------------------------------------
$wsRights = $wb.Worksheets.Item("RightsStatement")

$wsRights.Range("E2").Value = "Unless expressly stated otherwise, the organization that has made this Item available makes no warranties about the Item and cannot guarantee the accuracy of this Rights Statement. You are responsible for your own use."
$wsRights.Range("E4").Value = "Unless expressly stated otherwise, the organization that has made this Item available makes no warranties about the Item and cannot guarantee the accuracy of this Rights Statement. You are responsible for your own use."
$wsRights.Range("E5").Value = "You may need to obtain other permissions for your intended use. For example, other rights such as publicity, privacy or moral rights may limit how you may use the material."
$wsRights.Range("E6").Value = "You may need to obtain other permissions for your intended use. For example, other rights such as publicity, privacy or moral rights may limit how you may use the material."
$wsRights.Range("E7").Value = "You may find additional information about the copyright status of the Item on the website of the organization that has made the Item available."
$wsRights.Range("E8").Value = "You may find additional information about the copyright status of the Item on the website of the organization that has made the Item available."
$wsRights.Range("E10").Value = "You may need to obtain other permissions for your intended use. For example, other rights such as publicity, privacy, or moral rights may limit how you may use the material."
$wsRights.Range("E11").Value = "Unless expressly stated otherwise, the organization that has made this Item available makes no warranties about the Item and cannot guarantee the accuracy of this Rights Statement. You are responsible for your own use."
$wsRights.Range("E12").Value = "You may find additional information about the copyright status of the Item on the website of the organization that has made the Item available."
$wsRights.Range("E13").Value = "You may find additional information about the copyright status of the Item on the website of the organization that has made the Item available."

# ---------------------------------------------------------------------------
# 3. Image sheet: rightsHolderLiteral column (H) used to be a placeholder
#    "Property definition rights holder" for the Property-based images, and
#    missed the " image N" suffix for the Work/Institution based images.
#    It should now read "<dcterms:title> rights holder" for every image row.
# ---------------------------------------------------------------------------
$wsImage = $wb.Worksheets.Item("Image")
$lastRow = $wsImage.Cells.Item($wsImage.Rows.Count, 3).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $title = $wsImage.Cells.Item($r, 3).Text
    $wsImage.Cells.Item($r, 8).Value = "$title rights holder"
}

# ---------------------------------------------------------------------------
# 4. Person sheet: person4's page now points to Wikidata instead of Wikipedia
# ---------------------------------------------------------------------------
$wsPerson = $wb.Worksheets.Item("Person")
$wsPerson.Range("F6").Value = "http://www.wikidata.org/entity/Q7251"

# ---------------------------------------------------------------------------
# 5. Institution sheet: license, rights and rightsHolderLiteral columns were
#    dropped entirely, leaving just @id and name.
# ---------------------------------------------------------------------------
$wsInstitution = $wb.Worksheets.Item("Institution")
$wsInstitution.Columns.Item(5).Delete()
$wsInstitution.Columns.Item(4).Delete()
$wsInstitution.Columns.Item(2).Delete()
